$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 44305
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 22000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 22000
$ws.Range("S2").Value = 1100

$ws.Range("D3").Value = 44298
$ws.Range("M3").Value = 65
$ws.Range("N3").Value = 22000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 22000
$ws.Range("S3").Value = 1100

$ws.Range("D4").Value = 44382
$ws.Range("M4").Value = 24
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 20000
$ws.Range("S4").Value = 1000

$ws.Range("D5").Value = 44377
$ws.Range("M5").Value = 25
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("S5").Value = 1000

$ws.Range("D6").Value = 44403
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 20000
$ws.Range("S6").Value = 1000

$ws.Range("D7").Value = 44445
$ws.Range("M7").Value = 45
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 20000
$ws.Range("S7").Value = 1000

$ws.Range("D8").Value = 44389
$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 20000
$ws.Range("S8").Value = 1000

$ws.Range("D9").Value = 44385
$ws.Range("M9").Value = 36
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("S9").Value = 1000

$ws.Range("D10").Value = 44448
$ws.Range("M10").Value = 30
$ws.Range("N10").Value = 22000
$ws.Range("O10").Value = 22000
$ws.Range("P10").Value = 22000
$ws.Range("S10").Value = 1100

$ws.Range("D11").Value = 44307
$ws.Range("M11").Value = 30
$ws.Range("N11").Value = 22000
$ws.Range("O11").Value = 22000
$ws.Range("P11").Value = 22000
$ws.Range("S11").Value = 1100

$ws.Range("D12").Value = 44406
$ws.Range("M12").Value = 20
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 20000
$ws.Range("S12").Value = 1000

$ws.Range("D13").Value = 44294
$ws.Range("M13").Value = 25
$ws.Range("N13").Value = 25000
$ws.Range("O13").Value = 25000
$ws.Range("P13").Value = 25000
$ws.Range("S13").Value = 1250

$ws.Range("D14").Value = 44413
$ws.Range("M14").Value = 45
$ws.Range("N14").Value = 20000
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 20000
$ws.Range("S14").Value = 1000

$ws.Range("D15").Value = 44291
$ws.Range("M15").Value = 70
$ws.Range("N15").Value = 25000
$ws.Range("O15").Value = 25000
$ws.Range("P15").Value = 25000
$ws.Range("S15").Value = 1250

$ws.Range("D16").Value = 44292
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = 25000
$ws.Range("O16").Value = 25000
$ws.Range("P16").Value = 25000
$ws.Range("S16").Value = 1250

$ws.Range("D17").Value = 44300
$ws.Range("M17").Value = 45
$ws.Range("N17").Value = 22000
$ws.Range("O17").Value = 22000
$ws.Range("P17").Value = 22000
$ws.Range("S17").Value = 1100

$ws.Range("D18").Value = 44301
$ws.Range("M18").Value = 38
$ws.Range("N18").Value = 22000
$ws.Range("O18").Value = 22000
$ws.Range("P18").Value = 22000
$ws.Range("S18").Value = 1100

$ws.Range("D19").Value = 44400
$ws.Range("M19").Value = 45
$ws.Range("N19").Value = 20000
$ws.Range("O19").Value = 20000
$ws.Range("P19").Value = 20000
$ws.Range("S19").Value = 1000

$ws.Range("D20").Value = 44376
$ws.Range("M20").Value = 38
$ws.Range("N20").Value = 20000
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 20000
$ws.Range("S20").Value = 1000
